$wb = $excel.ActiveWorkbook

# Grab the two existing worksheets by tab position (left-to-right), BEFORE any renaming:
#   position 1 -> currently named "Sheet1" (physical sheet1.xml, not the active tab)
#   position 2 -> currently named "Sheet2" (physical sheet2.xml, the active/selected tab)
$firstSheet  = $wb.Worksheets.Item(1)
$secondSheet = $wb.Worksheets.Item(2)

$values = @("uipath", "hasini", "akira", "nandan")

# Replace the contents of both sheets with the same new data set (A1:A4),
# clearing out whatever used to be there (including now-unused shared strings).
foreach ($sh in @($firstSheet, $secondSheet)) {
    $sh.Cells.Clear()
    for ($i = 0; $i -lt $values.Length; $i++) {
        $sh.Cells.Item($i + 1, 1).Value = $values[$i]
    }
}

# Mark the selection on the sheet that will stay in the background (position 1).
[void]$firstSheet.Range("A1:A4").Select()

# Swap the two sheets' names (go through a temporary name to avoid a collision).
$firstSheet.Name = "__tmp_swap__"
$secondSheet.Name = "Sheet1"
$firstSheet.Name = "Sheet2"

# The sheet now named "Sheet1" (originally in position 2) is the one that stays active/selected.
[void]$secondSheet.Range("A1:A4").Select()
[void]$secondSheet.Activate()
